$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 206.83333
$ws.Range("I2").Value = 48.2
$ws.Range("K2").Value = 48.2
$ws.Range("M2").Value = 64.8
$ws.Range("H18").Value = 356.66666
$ws.Range("I18").Value = 356.66666
$ws.Range("K18").Value = 356.66666
$ws.Range("M18").Value = -72.66665999999998
$ws.Range("H40").Value = 1764.7142
$ws.Range("I40").Value = 1666.5
$ws.Range("J40").Value = 1838.375
$ws.Range("K40").Value = 1666.5
$ws.Range("L40").Value = 1838.375
$ws.Range("M40").Value = -1491.5
$ws.Range("N40").Value = -2188.375
$ws.Range("H53").Value = 198.05882
$ws.Range("I53").Value = 170.4
$ws.Range("J53").Value = 209.58333
$ws.Range("K53").Value = 170.4
$ws.Range("L53").Value = 209.58333
$ws.Range("M53").Value = 466.6
$ws.Range("N53").Value = -1483.58333
$ws.Range("H100").Value = 1876.875
$ws.Range("I100").Value = 987.5
$ws.Range("J100").Value = 2766.25
$ws.Range("K100").Value = 987.5
$ws.Range("L100").Value = 2766.25
$ws.Range("M100").Value = -446.5
$ws.Range("N100").Value = -3848.25
$ws.Range("H132").Value = 1039.4736
$ws.Range("I132").Value = 995.36365
$ws.Range("J132").Value = 2252.5
$ws.Range("K132").Value = 2986.09095
$ws.Range("L132").Value = 6757.5
$ws.Range("M132").Value = -456.0909499999998
$ws.Range("N132").Value = -11817.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2035.4546
$ws.Range("I97").Value = 1581.6666
$ws.Range("J97").Value = 2580
$ws.Range("K97").Value = 1581.6666
$ws.Range("L97").Value = 2580
$ws.Range("M97").Value = -1085.6666
$ws.Range("N97").Value = -3572
$ws.Range("H102").Value = 2205.1052
$ws.Range("I102").Value = 1914
$ws.Range("K102").Value = 1914
$ws.Range("M102").Value = -292
$ws.Range("H132").Value = 6879.593
$ws.Range("I132").Value = 2408.8462
$ws.Range("J132").Value = 11031
$ws.Range("K132").Value = 7226.5386
$ws.Range("L132").Value = 33093
$ws.Range("M132").Value = -4696.5386
$ws.Range("N132").Value = -38153

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2672.3333
$ws.Range("I20").Value = 1008
$ws.Range("J20").Value = 3504.5
$ws.Range("K20").Value = 1008
$ws.Range("L20").Value = 3504.5
$ws.Range("M20").Value = -761
$ws.Range("N20").Value = -3998.5
$ws.Range("H94").Value = 1124.1364
$ws.Range("I94").Value = 959
$ws.Range("J94").Value = 1341.421
$ws.Range("K94").Value = 959
$ws.Range("L94").Value = 1341.421
$ws.Range("M94").Value = -508
$ws.Range("N94").Value = -2243.421
$ws.Range("H99").Value = 2226.4736
$ws.Range("I99").Value = 2191
$ws.Range("J99").Value = 2265.889
$ws.Range("K99").Value = 2191
$ws.Range("L99").Value = 2265.889
$ws.Range("M99").Value = -693
$ws.Range("N99").Value = -5261.889
$ws.Range("H105").Value = 3509.1707
$ws.Range("I105").Value = 2914.6072
$ws.Range("K105").Value = 2914.6072
$ws.Range("M105").Value = -1167.6072
$ws.Range("H107").Value = 1679.5333
$ws.Range("I107").Value = 1666.3636
$ws.Range("J107").Value = 1715.75
$ws.Range("K107").Value = 1666.3636
$ws.Range("L107").Value = 1715.75
$ws.Range("M107").Value = 253.6364000000001
$ws.Range("N107").Value = -5555.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 103.888885
$ws.Range("I22").Value = 72.5
$ws.Range("J22").Value = 166.66667
$ws.Range("K22").Value = 72.5
$ws.Range("L22").Value = 166.66667
$ws.Range("M22").Value = 277.5
$ws.Range("N22").Value = -866.6666700000001
$ws.Range("H127").Value = 15000000
$ws.Range("J127").Value = 15000000
$ws.Range("L127").Value = 15000000
$ws.Range("N127").Value = -15009920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3841.0588
$ws.Range("I123").Value = 1466.6666
$ws.Range("J123").Value = 4349.857
$ws.Range("K123").Value = 4399.9998
$ws.Range("L123").Value = 13049.571
$ws.Range("M123").Value = -1949.9998
$ws.Range("N123").Value = -17949.571
$ws.Range("H125").Value = 2782.6924
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 2847.9167
$ws.Range("K125").Value = 6000
$ws.Range("L125").Value = 8543.750100000001
$ws.Range("M125").Value = -1080
$ws.Range("N125").Value = -18383.7501
$ws.Range("H132").Value = 1923.6316
$ws.Range("I132").Value = 3809.6667
$ws.Range("J132").Value = 1570
$ws.Range("K132").Value = 34287.0003
$ws.Range("L132").Value = 14130
$ws.Range("M132").Value = -31757.0003
$ws.Range("N132").Value = -19190

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5505.5557
$ws.Range("I80").Value = 8985.714
$ws.Range("J80").Value = 3290.9092
$ws.Range("K80").Value = 8985.714
$ws.Range("L80").Value = 3290.9092
$ws.Range("M80").Value = -7987.714
$ws.Range("N80").Value = -5286.9092
$ws.Range("H83").Value = 5505.5557
$ws.Range("I83").Value = 8985.714
$ws.Range("J83").Value = 3290.9092
$ws.Range("K83").Value = 44928.57
$ws.Range("L83").Value = 16454.546
$ws.Range("M83").Value = -39936.57
$ws.Range("N83").Value = -26438.546
$ws.Range("H97").Value = 2750
$ws.Range("I97").Value = 2400
$ws.Range("J97").Value = 3100
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 3100
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -4092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 327.13333
$ws.Range("I22").Value = 240
$ws.Range("J22").Value = 426.7143
$ws.Range("K22").Value = 240
$ws.Range("L22").Value = 426.7143
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = -1016.7143
$ws.Range("H27").Value = 327.13333
$ws.Range("I27").Value = 240
$ws.Range("J27").Value = 426.7143
$ws.Range("K27").Value = 240
$ws.Range("L27").Value = 426.7143
$ws.Range("M27").Value = -133
$ws.Range("N27").Value = -640.7143
$ws.Range("H46").Value = 1099.8334
$ws.Range("I46").Value = 999.75
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 999.75
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -811.75
$ws.Range("N46").Value = -1676
$ws.Range("H93").Value = 777.6667
$ws.Range("I93").Value = 777.6667
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 777.6667
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 470.3333
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2171.5715
$ws.Range("I81").Value = 1433.6666
$ws.Range("J81").Value = 2725
$ws.Range("K81").Value = 2867.3332
$ws.Range("L81").Value = 5450
$ws.Range("M81").Value = -1806.3332
$ws.Range("N81").Value = -7572
$ws.Range("H84").Value = 2171.5715
$ws.Range("I84").Value = 1433.6666
$ws.Range("J84").Value = 2725
$ws.Range("K84").Value = 14336.666
$ws.Range("L84").Value = 27250
$ws.Range("M84").Value = -9032.666000000001
$ws.Range("N84").Value = -37858
$ws.Range("H107").Value = 3674.875
$ws.Range("I107").Value = 1679
$ws.Range("K107").Value = 5037
$ws.Range("M107").Value = -3117
